$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # "loginTest"
$ws2 = $wb.Worksheets.Item(2)   # "Sheet1"

# --- Rebuild "Sheet1" (ws2) content -------------------------------------
# Remove existing hyperlinks and wipe the old data entirely.
$ws2.Hyperlinks.Delete()
$ws2.Cells.Clear()

$ws2.Range("A1").Value = "username"
$ws2.Range("B1").Value = "password"
$ws2.Range("C1").Value = "status"

$ws2.Range("A2").Value = "wrong_username"
$ws2.Range("B2").Value = "wrong_password"
$ws2.Range("C2").Value = "failed"

$ws2.Range("A3").Value = "wrong_username"
$ws2.Range("B3").Value = "wrong_password"
$ws2.Range("C3").Value = "failed"

$ws2.Range("A4").Value = "wrong_username"
$ws2.Range("B4").Value = "wrong_password"
$ws2.Range("C4").Value = "failed"

$ws2.Range("A5").Value = "correct_username"
$ws2.Range("B5").Value = "correct_password"
$ws2.Range("C5").Value = "passed"

$ws2.Columns.Item(1).ColumnWidth = 17.6640625
$ws2.Columns.Item(2).ColumnWidth = 17.5546875

# --- Adjust selections ---------------------------------------------------
$ws1.Range("A1:B5").Select()
$ws2.Range("C5").Select()

# Sheet1 ("Sheet1" tab) ends up as the active / selected tab.
$ws2.Activate()
$ws2.Range("C5").Select()

# --- Remove now-unused "Hyperlink" cell style ----------------------------
$wb.Styles.Item("Hyperlink").Delete()
